# Updated symbol list on Mon Jan 30 11:26:44 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns for the crypto rows that moved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $Text into $Range as literal text (matching the source data's
# inline-string cells) instead of letting Excel auto-convert numeric-looking
# strings (e.g. "308.34", "-2.09%") into Number/Percent cells. Forcing the
# "@" text format for the assignment, then clearing formats afterwards, keeps
# the value as text while leaving the cell style untouched (General/no style).
function Set-TextValue($Range, $Text) {
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "308.34"
Set-TextValue $ws.Range("E2") "-2.09%"
Set-TextValue $ws.Range("D3") "37.77"
Set-TextValue $ws.Range("E3") "-4.16%"
Set-TextValue $ws.Range("E4") "-2.05%"
Set-TextValue $ws.Range("D5") "0.07895"
Set-TextValue $ws.Range("E5") "-3.37%"
Set-TextValue $ws.Range("D6") "2.022"
Set-TextValue $ws.Range("E6") "2.66%"
Set-TextValue $ws.Range("D7") "4.392"
Set-TextValue $ws.Range("E7") "3.47%"
Set-TextValue $ws.Range("D8") "8.225"
Set-TextValue $ws.Range("E8") "0.09%"
Set-TextValue $ws.Range("E9") "-3.00%"
Set-TextValue $ws.Range("D10") "0.9267"
Set-TextValue $ws.Range("E10") "0.18%"
Set-TextValue $ws.Range("D11") "0.1280"
Set-TextValue $ws.Range("E11") "-8.94%"
Set-TextValue $ws.Range("D12") "0.1904"
Set-TextValue $ws.Range("E12") "-3.60%"
Set-TextValue $ws.Range("D13") "0.08798"
Set-TextValue $ws.Range("E13") "-2.43%"
Set-TextValue $ws.Range("D14") "0.03454"
Set-TextValue $ws.Range("E14") "-1.40%"
Set-TextValue $ws.Range("D15") "0.09748"
Set-TextValue $ws.Range("E15") "-0.84%"
Set-TextValue $ws.Range("E16") "-0.75%"
Set-TextValue $ws.Range("D17") "0.006164"
Set-TextValue $ws.Range("E17") "2.24%"
Set-TextValue $ws.Range("D18") "3.548"
Set-TextValue $ws.Range("E18") "-2.93%"
Set-TextValue $ws.Range("D19") "0.3441"
Set-TextValue $ws.Range("E19") "-0.52%"
Set-TextValue $ws.Range("E20") "-3.30%"
Set-TextValue $ws.Range("D21") "5.009"
Set-TextValue $ws.Range("E21") "5.29%"
Set-TextValue $ws.Range("D22") "0.2518"
Set-TextValue $ws.Range("E22") "3.92%"
Set-TextValue $ws.Range("D23") "0.04328"
Set-TextValue $ws.Range("E23") "-0.88%"
Set-TextValue $ws.Range("D24") "0.001226"
Set-TextValue $ws.Range("E24") "0.29%"
Set-TextValue $ws.Range("D25") "0.004594"
Set-TextValue $ws.Range("E26") "176.73%"
Set-TextValue $ws.Range("D39") "0.02255"
Set-TextValue $ws.Range("E39") "4.55%"
Set-TextValue $ws.Range("D40") "0.05029"
Set-TextValue $ws.Range("E40") "-3.27%"
Set-TextValue $ws.Range("D41") "0.007540"
Set-TextValue $ws.Range("E41") "0.07%"
Set-TextValue $ws.Range("D42") "0.009930"
Set-TextValue $ws.Range("E42") "0.77%"
Set-TextValue $ws.Range("E43") "-1.19%"
Set-TextValue $ws.Range("D44") "0.002094"
Set-TextValue $ws.Range("E44") "-0.91%"
Set-TextValue $ws.Range("D45") "0.008551"
Set-TextValue $ws.Range("E45") "-6.24%"
Set-TextValue $ws.Range("D46") "0.00006439"
Set-TextValue $ws.Range("E46") "0.52%"
Set-TextValue $ws.Range("E47") "0.26%"
Set-TextValue $ws.Range("D48") "0.003006"
Set-TextValue $ws.Range("E48") "9.00%"
Set-TextValue $ws.Range("E49") "0.26%"
Set-TextValue $ws.Range("D50") "0.00002106"
Set-TextValue $ws.Range("E50") "0.26%"
Set-TextValue $ws.Range("D51") "0.0002006"
Set-TextValue $ws.Range("E51") "0.26%"
